# Language table update: add "victory" screen strings (victory/population/
# houses_deployed key-value pairs) and reword the "new_house" message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows right after the "new_house" row (row 19) to make
# room for the new victory/population/houses_deployed key-value pairs.
$ws.Rows("20:22").Insert()

# Fill in the new rows first so the shared-string table picks these values
# up before the reworded "new_house" text is written.
$ws.Cells.Item(20, 1).Value2 = "victory"
$ws.Cells.Item(20, 2).Value2 = "VICTORY"
$ws.Cells.Item(21, 1).Value2 = "population"
$ws.Cells.Item(21, 2).Value2 = "Population"
$ws.Cells.Item(22, 1).Value2 = "houses_deployed"
$ws.Cells.Item(22, 2).Value2 = "Homes Deployed"

# Reword the existing "new_house" message.
$ws.Cells.Item(19, 2).Value2 = "A new house is available! Deploy it to progress."

# Update the selected cell to match the edited workbook.
[void]$ws.Range("B21").Select()
